$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 12987415
$ws.Range("I19").Value = 28571744
$ws.Range("J19").Value = 475
$ws.Range("K19").Value = 28571744
$ws.Range("L19").Value = 475
$ws.Range("M19").Value = -28571569
$ws.Range("N19").Value = -825

# Row 33
$ws.Range("H33").Value = 291.7069
$ws.Range("I33").Value = 190.55319
$ws.Range("J33").Value = 723.9091
$ws.Range("K33").Value = 190.55319
$ws.Range("L33").Value = 723.9091
$ws.Range("M33").Value = 38.44681
$ws.Range("N33").Value = -1181.9091

# Row 98
$ws.Range("H98").Value = 1664.2413
$ws.Range("I98").Value = 1344.6666
$ws.Range("J98").Value = 3198.2
$ws.Range("K98").Value = 1344.6666
$ws.Range("L98").Value = 3198.2
$ws.Range("M98").Value = 153.3334
$ws.Range("N98").Value = -6194.2

# Row 107
$ws.Range("H107").Value = 872.3333
$ws.Range("I107").Value = 906.8421
$ws.Range("K107").Value = 906.8421
$ws.Range("M107").Value = 1013.1579

# Row 121
$ws.Range("H121").Value = 808
$ws.Range("I121").Value = 550
$ws.Range("J121").Value = 1324
$ws.Range("K121").Value = 1650
$ws.Range("L121").Value = 3972
$ws.Range("M121").Value = 97
$ws.Range("N121").Value = -7466

# Row 122
$ws.Range("H122").Value = 1664.2413
$ws.Range("I122").Value = 1344.6666
$ws.Range("J122").Value = 3198.2
$ws.Range("K122").Value = 4033.9998
$ws.Range("L122").Value = 9594.599999999999
$ws.Range("M122").Value = -1583.9998
$ws.Range("N122").Value = -14494.6

# Row 132
$ws.Range("H132").Value = 4257173.5
$ws.Range("I132").Value = 5264698
$ws.Range("J132").Value = 3180.2222
$ws.Range("K132").Value = 15794094
$ws.Range("L132").Value = 9540.6666
$ws.Range("M132").Value = -15791564
$ws.Range("N132").Value = -14600.6666

# Row 137
$ws.Range("H137").Value = 2567575
$ws.Range("I137").Value = 3229286.8
$ws.Range("J137").Value = 3442.875
$ws.Range("K137").Value = 9687860.399999999
$ws.Range("L137").Value = 10328.625
$ws.Range("M137").Value = -9685310.399999999
$ws.Range("N137").Value = -15428.625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6862.231
$ws.Range("I32").Value = 5922.419
$ws.Range("K32").Value = 5922.419
$ws.Range("M32").Value = -5635.419

# Row 63
$ws.Range("H63").Value = 2576.923
$ws.Range("I63").Value = 2136.3635
$ws.Range("K63").Value = 2136.3635
$ws.Range("M63").Value = -1450.3635

# Row 66
$ws.Range("H66").Value = 2576.923
$ws.Range("I66").Value = 2136.3635
$ws.Range("K66").Value = 10681.8175
$ws.Range("M66").Value = -7249.817499999999

# Row 102
$ws.Range("H102").Value = 2659.9524
$ws.Range("I102").Value = 2303.2778
$ws.Range("K102").Value = 2303.2778
$ws.Range("M102").Value = -681.2777999999998

# Row 132
$ws.Range("H132").Value = 3637.111
$ws.Range("I132").Value = 3281.3635
$ws.Range("J132").Value = 5202.4
$ws.Range("K132").Value = 9844.0905
$ws.Range("L132").Value = 15607.2
$ws.Range("M132").Value = -7314.0905
$ws.Range("N132").Value = -20667.2

$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 3728.5625
$ws.Range("I134").Value = 2741.7693
$ws.Range("J134").Value = 8004.6665
$ws.Range("K134").Value = 8225.3079
$ws.Range("L134").Value = 24013.9995
$ws.Range("M134").Value = -5690.3079
$ws.Range("N134").Value = -29083.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 2634.739
$ws.Range("I105").Value = 2224.95
$ws.Range("J105").Value = 5366.6665
$ws.Range("K105").Value = 2224.95
$ws.Range("L105").Value = 5366.6665
$ws.Range("M105").Value = -477.9499999999998
$ws.Range("N105").Value = -8860.666499999999

# Row 132
$ws.Range("H132").Value = 2221.8462
$ws.Range("I132").Value = 1759
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 5277
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -2747
$ws.Range("N132").Value = -15260

# Row 134
$ws.Range("H134").Value = 2040.7142
$ws.Range("I134").Value = 772.1875
$ws.Range("J134").Value = 6100
$ws.Range("K134").Value = 2316.5625
$ws.Range("L134").Value = 18300
$ws.Range("M134").Value = 218.4375
$ws.Range("N134").Value = -23370

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# Row 138
$ws.Range("H138").Value = 28615
$ws.Range("J138").Value = 28615
$ws.Range("L138").Value = 28615
$ws.Range("N138").Value = -38895

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1434.8667
$ws.Range("I107").Value = 448
$ws.Range("J107").Value = 1928.3
$ws.Range("K107").Value = 1344
$ws.Range("L107").Value = 5784.9
$ws.Range("M107").Value = 576
$ws.Range("N107").Value = -9624.9

# Row 131
$ws.Range("H131").Value = 1399.3704
$ws.Range("J131").Value = 1111.6888
$ws.Range("L131").Value = 3335.0664
$ws.Range("N131").Value = -13415.0664

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 20859.143
$ws.Range("J24").Value = 20859.143
$ws.Range("L24").Value = 20859.143
$ws.Range("N24").Value = -21205.143

# Row 70
$ws.Range("H70").Value = 4550.852
$ws.Range("I70").Value = 4597.125
$ws.Range("J70").Value = 4483.5454
$ws.Range("K70").Value = 4597.125
$ws.Range("L70").Value = 4483.5454
$ws.Range("M70").Value = -4327.125
$ws.Range("N70").Value = -5023.5454

# Row 73
$ws.Range("H73").Value = 4550.852
$ws.Range("I73").Value = 4597.125
$ws.Range("J73").Value = 4483.5454
$ws.Range("K73").Value = 4597.125
$ws.Range("L73").Value = 4483.5454
$ws.Range("M73").Value = -3661.125
$ws.Range("N73").Value = -6355.5454

# Row 102
$ws.Range("H102").Value = 38298.93
$ws.Range("I102").Value = 2039.1111
$ws.Range("J102").Value = 103566.6
$ws.Range("K102").Value = 2039.1111
$ws.Range("L102").Value = 103566.6
$ws.Range("M102").Value = -417.1111000000001
$ws.Range("N102").Value = -106810.6

# Row 113
$ws.Range("H113").Value = 1931.9048
$ws.Range("I113").Value = 1438
$ws.Range("J113").Value = 3166.6667
$ws.Range("K113").Value = 1438
$ws.Range("L113").Value = 3166.6667
$ws.Range("M113").Value = 732
$ws.Range("N113").Value = -7506.6667

# Row 132
$ws.Range("H132").Value = 4339.423
$ws.Range("I132").Value = 4701.9287
$ws.Range("J132").Value = 3916.5
$ws.Range("K132").Value = 14105.7861
$ws.Range("L132").Value = 11749.5
$ws.Range("M132").Value = -11575.7861
$ws.Range("N132").Value = -16809.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8335474
$ws.Range("I7").Value = 25000872
$ws.Range("J7").Value = 2775
$ws.Range("K7").Value = 25000872
$ws.Range("L7").Value = 2775
$ws.Range("M7").Value = -25000760
$ws.Range("N7").Value = -2999

# Row 46
$ws.Range("H46").Value = 1507.5807
$ws.Range("I46").Value = 953.46155
$ws.Range("J46").Value = 4389
$ws.Range("K46").Value = 953.46155
$ws.Range("L46").Value = 4389
$ws.Range("M46").Value = -765.46155
$ws.Range("N46").Value = -4765

# Row 122
$ws.Range("H122").Value = 3499.9285
$ws.Range("I122").Value = 2624.875
$ws.Range("K122").Value = 7874.625
$ws.Range("M122").Value = -5424.625

# Row 126
$ws.Range("H126").Value = 8335474
$ws.Range("I126").Value = 25000872
$ws.Range("J126").Value = 2775
$ws.Range("K126").Value = 75002616
$ws.Range("L126").Value = 8325
$ws.Range("M126").Value = -75000146
$ws.Range("N126").Value = -13265

# Row 132
$ws.Range("H132").Value = 2544.9756
$ws.Range("I132").Value = 1937.3
$ws.Range("J132").Value = 3123.7144
$ws.Range("K132").Value = 5811.9
$ws.Range("L132").Value = 9371.143199999999
$ws.Range("M132").Value = -3281.9
$ws.Range("N132").Value = -14431.1432

# Row 136
$ws.Range("H136").Value = 5886123.5
$ws.Range("I136").Value = 7693700
$ws.Range("J136").Value = 11500
$ws.Range("K136").Value = 23081100
$ws.Range("L136").Value = 34500
$ws.Range("M136").Value = -23078550
$ws.Range("N136").Value = -39600

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()

# Row 140
$ws.Range("H140").Value = 29600
$ws.Range("J140").Value = 29600
$ws.Range("L140").Value = 29600
$ws.Range("N140").Value = -39960

$ws = $wb.Worksheets.Item("WVR")
# Row 40
$ws.Range("H40").Value = 26590.334
$ws.Range("J40").Value = 26590.334
$ws.Range("L40").Value = 26590.334
$ws.Range("N40").Value = -26888.334

# Row 74
$ws.Range("H74").Value = 11942.571
$ws.Range("J74").Value = 11942.571
$ws.Range("L74").Value = 11942.571
$ws.Range("N74").Value = -13814.571

# Row 77
$ws.Range("H77").Value = 11942.571
$ws.Range("J77").Value = 11942.571
$ws.Range("L77").Value = 35827.713
$ws.Range("N77").Value = -45187.713

# Row 136
$ws.Range("H136").Value = 1560.8928
$ws.Range("I136").Value = 754.1177
$ws.Range("J136").Value = 2807.7273
$ws.Range("K136").Value = 2262.3531
$ws.Range("L136").Value = 8423.1819
$ws.Range("M136").Value = 287.6468999999997
$ws.Range("N136").Value = -13523.1819

Write-Output "Applied market-price refresh across all sheets."
